# Apply attendance updates for week 4 (column F) and update the selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the attendance table that get a "1" (present) for săpt. 4 (column F)
$rows = @(6, 9, 10, 12, 13, 14, 17, 18, 19, 21)

foreach ($r in $rows) {
    $ws.Range("F$r").Value = 1
}

# Update the active selection on the frozen (bottom-right) pane to F3:F21,
# with F3 as the active cell.
$ws.Range("F3:F21").Select()
